$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 36 (pushes existing rows 36-44 down to 37-45,
# and copies formatting down from row 35 into the new row 36, matching
# Excel's native "Insert" behaviour).
$ws.Rows("36").Insert()

# Populate the new row 36 with its label and row height; clear any
# carried-over content in W36 (format-only, no value).
$ws.Range("B36").Value = "Validierungsdatensatz"
$ws.Rows("36").RowHeight = 28.8
$ws.Range("W36").Value = $null

# Update sheet view to match the saved state after the edit.
$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Range("D37").Select() | Out-Null
